$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Desk Research" team-member time value (C9) and
# bump the "Start paper prototype" time value (C12) from 5 to 6.
$ws.Range("C9").Value = 4
$ws.Range("C12").Value = 6

# Move the active selection from C18 to C19, matching the saved cursor
# position in the workbook.
$ws.Range("C19").Select()
